# Beavers incremental mega sheet - add "Sheet1" (levelling / pacing calculations)
# after the "Fishing" sheet, and make it the active sheet.

$wb = $excel.ActiveWorkbook

# --- Add the new worksheet after the last (Fishing) sheet -------------------
$fishing = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $fishing)
# Excel names a freshly added sheet "Sheet1" by default - that matches the
# target workbook, so nothing else to rename.

# --- Header / explanation cells (row 3 & 4) ---------------------------------
# Written in this exact order so the shared-strings table gets the same
# index assignment as the target file (907..911).
$ws.Range("D3").Value = "round(pow(1.6, level) + 1.5 * pow(level, 4.3)) + 1 * (level + 1) * 7 + 2"
$ws.Range("C4").Value = "level"
$ws.Range("D4").Value = "2^C5 + 2 * (C5 + 1) * 7 + 2"
$ws.Range("F4").Value = "2^C5 + 2 * (C5^3) + 1 * (C5 + 1) * 7 + 2"
$ws.Range("E4").Value = "3^C5 + 2 * (C5^2) + 1 * (C5 + 1) * 7 + 2"

# --- Row 5: first data row, plain (non-shared) formulas ---------------------
$ws.Range("C5").Value = 1
$ws.Range("D5").Formula = "=2^C5 + 2 * (C5^2) + 1 * (C5+ 1) * 7 + 2"
$ws.Range("E5").Formula = "=3^C5 + 2 * (C5^2) + 1 * (C5+ 1) * 7 + 2"
$ws.Range("F5").Formula = "=2^C5 + 2 * (C5^3) + 1 * (C5 + 1) * 7 + 2"

# --- Rows 6-46: level numbers + shared formulas ------------------------------
for ($r = 6; $r -le 46; $r++) {
    $ws.Range("C$r").Value = $r - 4
}
$ws.Range("D6:D46").Formula = "=2^C6 + 2 * (C6^2) + 1 * (C6+ 1) * 7 + 2"
$ws.Range("E6:E46").Formula = "=3^C6 + 2 * (C6^2) + 1 * (C6+ 1) * 7 + 2"
$ws.Range("F6:F46").Formula = "=2^C6 + 2 * (C6^3) + 1 * (C6 + 1) * 7 + 2"

# --- Column widths (D=57, E=36, F=~35.86 characters) -------------------------
$ws.Columns.Item(4).ColumnWidth = 56.166666666666664
$ws.Columns.Item(5).ColumnWidth = 35.166666666666664
$ws.Columns.Item(6).ColumnWidth = 35.02134

# --- Selection / active cell on the new sheet --------------------------------
$ws.Range("H63").Select()

# The sheet is now the last one and becomes the active tab, matching the
# workbook's activeTab/bookViews update in the target file.
